$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-05 Wednesday" "2024-06-06 Thursday"
Replace-Text "24×42=" "47×54="
Replace-Text "93×39=" "66×89="
Replace-Text "89×18=" "37×16="
Replace-Text "34×37=" "94×87="
Replace-Text "18×52=" "48×25="
Replace-Text "27×27=" "23×61="
Replace-Text "85×54=" "46×42="
Replace-Text "98×95=" "36×38="
Replace-Text "92×35=" "27×11="
Replace-Text "96×84=" "82×50="
Replace-Text "84×51=" "81×63="
Replace-Text "87×57=" "74×27="
Replace-Text "73×81=" "55×46="
Replace-Text "24×30=" "79×27="
Replace-Text "27×67=" "29×98="
Replace-Text "90×50=" "11×53="
Replace-Text "60×88=" "82×78="
Replace-Text "52×37=" "71×76="
Replace-Text "73×75=" "38×17="
Replace-Text "21×81=" "72×47="
Replace-Text "65×20=" "85×28="
Replace-Text "34×25=" "77×90="
Replace-Text "85×75=" "11×67="
Replace-Text "59×19=" "31×23="
Replace-Text "53×33=" "49×52="
